# Update visitor-count figures (column F) on the "展览" and "全部类型" sheets
# to reflect newly generated output data (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 103
$ws1.Range("F4").Value = 7360
$ws1.Range("F6").Value = 438
$ws1.Range("F7").Value = 3883
$ws1.Range("F9").Value = 553
$ws1.Range("F11").Value = 619
$ws1.Range("F12").Value = 112

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 103
$ws4.Range("F5").Value = 7360
$ws4.Range("F8").Value = 438
$ws4.Range("F9").Value = 3883
$ws4.Range("F11").Value = 553
$ws4.Range("F13").Value = 619
$ws4.Range("F14").Value = 112
